# investors.xlsx -- "Added pivoting in filters"
#
# Row-1 header cells gain a short annotation (mandatory/optional
# asterisks moved to the fields that are actually required) and each
# header cell gets an explanatory cell comment describing the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook was authored under the name "Author"; match it so any
# comment/author metadata the runtime does expose lines up.
$excel.UserName = "Author"

# --- Row 1 header text -------------------------------------------------
$ws.Range("B1").Value = "PAN"
$ws.Range("C1").Value = "Primary Email *"
$ws.Range("D1").Value = "Tags"
$ws.Range("E1").Value = "Category *"

# --- Header-cell comments -----------------------------------------------
$ws.Range("A1").ClearComments()
[void]$ws.Range("A1").AddComment("Author:`n-Mandatory`n-Name you wish to record for finding investor`n- Name needs to be unique for your entity.  `n- Investing entity name can be different, which you can add in KYC")
$ws.Range("B1").ClearComments()
[void]$ws.Range("B1").AddComment("Author:`nOptional to add Tax ID No.  `nIn case available, helps combine investors added with different names`n")
$ws.Range("C1").ClearComments()
[void]$ws.Range("C1").AddComment("Author:`nMandatory`nThis is the unique identifier for each investor.  `nNote – this will not trigger any notification / access.  `n")
$ws.Range("D1").ClearComments()
[void]$ws.Range("D1").AddComment("Author:`nOptional`nTo help you find investor later  Can be whatever you want it to be! `n")
$ws.Range("E1").ClearComments()
[void]$ws.Range("E1").AddComment("Author:`nMandatory`nThis is to group investors in category you define.  `nCan edit Categories under the Home button / Entity details (home icon on top right)`n")
$ws.Range("F1").ClearComments()
[void]$ws.Range("F1").AddComment("Author:`nOptional`nTo help you find investors later.  Can be whatever you want it to be")

# --- Saved cursor/selection ---------------------------------------------
[void]$ws.Range("A7").Select()

